# PM04 Tidsregistrering for Bille.xlsx
# Fills in the previously-empty time-tracking rows 20-23 on sheet "Ark1"
# with task description, role, date, start time and end time, matching
# the author's new entries. Formulas in columns G/H already exist on the
# sheet and recalculate automatically once D/E/C receive values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# --- Formatting -----------------------------------------------------
# Row 19 already carries the exact cell formats (date / time / borders)
# that rows 20-23 need for columns A and C:H, so reuse it as a template.
$ws.Range("A19").Copy()
$ws.Range("A20:A23").PasteSpecial(-4122)

$ws.Range("C19:H19").Copy()
$ws.Range("C20:H20").PasteSpecial(-4122)
$ws.Range("C19:H19").Copy()
$ws.Range("C21:H21").PasteSpecial(-4122)
$ws.Range("C19:H19").Copy()
$ws.Range("C22:H22").PasteSpecial(-4122)
$ws.Range("C19:H19").Copy()
$ws.Range("C23:H23").PasteSpecial(-4122)

# Column B (Rolle) needs its own style - borrow it from a row that has
# a role entered already (row 3).
$ws.Range("B3").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B23").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Values -----------------------------------------------------------
# Task descriptions / roles (order chosen so new shared-string entries
# line up with the rows/cells that reference them).
$ws.Range("A20").Value = "review af DCD SD0104 + før forbedrelse pga. troede jeg skulle lave den"
$ws.Range("B20").Value = "Reviewer"

$ws.Range("A21").Value = "spilder tid pga. folk ikke skrev sig på opgaven jeg skrev mig på"

$ws.Range("A23").Value = "UT05, UT06"
$ws.Range("B23").Value = "Test Desinger"

$ws.Range("A22").Value = "UTD05 UTD05"
$ws.Range("B22").Value = "Test Analyst"

# Dates / start / end times
$ws.Range("C20").Value = 43892
$ws.Range("D20").Value = 0.3611111111111111
$ws.Range("E20").Value = 0.40972222222222227

$ws.Range("C21").Value = 43892
$ws.Range("D21").Value = 0.42708333333333331
$ws.Range("E21").Value = 0.47916666666666669

$ws.Range("C22").Value = 43893
$ws.Range("D22").Value = 0.375
$ws.Range("E22").Value = 0.46875

$ws.Range("C23").Value = 43893
$ws.Range("D23").Value = 0.46875
$ws.Range("E23").Value = 0.59375

# --- Selection ----------------------------------------------------------
$ws.Activate()
$ws.Range("B24").Select()
